# update int to long
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3 gets a numeric value
$ws.Range("A3").Value = 9999

# A4 gets a text value "a" (stored as shared string), with style matching
# the "right aligned / bordered" look already used in rows 10-11 (style 6)
$ws.Range("A4").Borders.LineStyle = 1
$ws.Range("A4").HorizontalAlignment = -4152
$ws.Range("A4").Value = "a"

# Update the visible selection to C6
$ws.Range("C6").Select()
